# Implement "Trash" feature rows in the FileUploadController section (soft
# delete of files) + a new "Bulk Download File" row, as per the commit:
# "Implemented Trash features with soft delete of files[ADD]"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new rows at the right spots -------------------------------
# 1 new row before the existing row 9 ("Get All Versions of a File" row,
# which becomes row 10) for the new "Bulk Download File" entry.
$ws.Rows.Item(9).Insert()

# 4 new rows before what is now row 13 (the "UserController" / "Get All
# Users" row) for the new Trash feature block.
$ws.Range("A13:A16").EntireRow.Insert()

# --- Row 9: Bulk Download File ---------------------------------------------
$ws.Cells.Item(9,1).Value = "FileUploadController"
$ws.Cells.Item(9,2).Value = "Bulk Download File"
$ws.Cells.Item(9,3).Value = "POST"
$ws.Cells.Item(9,4).Value = "/api/v1/files/bulk-download"
$ws.Cells.Item(9,5).Value = "Downloads file as ZIP format"
$ws.Cells.Item(9,6).Value = "Admin, Staff"
$ws.Cells.Item(9,7).Value = "{`n  `"fileIds`": [`n    `"3fa85f64-5717-4562-b3fc-2c963f66afa6`"`n  ],`n  `"zipFileName`": `"string`"`n}"
$ws.Cells.Item(9,8).Value = "Binary file stream"
$ws.Rows.Item(9).RowHeight = 96

# --- Row 13: Get all files in trash ----------------------------------------
$ws.Cells.Item(13,1).Value = "FileUploadController"
$ws.Cells.Item(13,2).Value = "Get all files in trash"
$ws.Cells.Item(13,3).Value = "GET"
$ws.Cells.Item(13,4).Value = "/api/v1/files/trash"
$ws.Cells.Item(13,5).Value = "Returns list of files details in the trash"
$ws.Cells.Item(13,6).Value = "Admin"
$ws.Cells.Item(13,7).Value = "None"
$ws.Rows.Item(13).RowHeight = 16

# --- Row 14: Restore from trash --------------------------------------------
$ws.Cells.Item(14,1).Value = "FileUploadController"
$ws.Cells.Item(14,2).Value = "Restore from trash"
$ws.Cells.Item(14,3).Value = "POST"
$ws.Cells.Item(14,4).Value = "/api/v1/files/trash/{id}/restore"
$ws.Cells.Item(14,5).Value = "Move file form trash to Archive file & make db update"
$ws.Cells.Item(14,6).Value = "Admin"
$ws.Cells.Item(14,7).Value = "id"
$ws.Rows.Item(14).RowHeight = 32

# --- Row 15: Delete permanently in trash ------------------------------------
$ws.Cells.Item(15,1).Value = "FileUploadController"
$ws.Cells.Item(15,2).Value = "Delete permanently in trash"
$ws.Cells.Item(15,3).Value = "DELETE"
$ws.Cells.Item(15,4).Value = "/api/v1/files/trash/{id}/permanent"
$ws.Cells.Item(15,5).Value = "Removes file at the trash & also in archieve if exist"
$ws.Cells.Item(15,6).Value = "Admin"
$ws.Cells.Item(15,7).Value = "id"
$ws.Rows.Item(15).RowHeight = 32

# --- Row 16: Get trash stats -------------------------------------------------
$ws.Cells.Item(16,1).Value = "FileUploadController"
$ws.Cells.Item(16,2).Value = "Get trash stats"
$ws.Cells.Item(16,3).Value = "GET"
$ws.Cells.Item(16,4).Value = "/api/v1/files/trash/stats"
$ws.Cells.Item(16,5).Value = "Get trash statistics"
$ws.Cells.Item(16,6).Value = "Admin"
$ws.Cells.Item(16,7).Value = "None"
$ws.Rows.Item(16).RowHeight = 16

# --- View state: scroll + selection, matching the saved workbook view ------
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$null = $ws.Range("E15").Select()
